# "Updated Elixir example usages"
# Slide 4 ("Enter elixir") body placeholder: update the bullet that lists
# which companies use Elixir, adding Spotify to the list.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item(2)          # "Text Placeholder 4"
$tf = $shp.TextFrame
$tr = $tf.TextRange

$para = $tr.Paragraphs(4)         # "Used by Discord and Pinterest"
$run = $para.Runs(1)
$run.Text = "Used by Discord, Pinterest, Spotify"
